# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The account-statement table (rows 16-21) is re-keyed to a new set of
# overdue periods (2506/2507/2508 instead of 2505/2506/2507), and the two
# workers' rows are now interleaved by period instead of grouped by worker.
# The "VALOR MORA" total (E11) is recalculated to match the new totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for the worker/period table (columns B..G, rows 16-21)
# B = doc type, C = doc number, D = worker name, E = period, F = monthly
# value, G = base salary.
$data = @(
    @("CC", "45523042", "NOHORA MILENA BERRIO CRECIAN", "2506", 80000, 2000000),
    @("CC", "73209840", "ALVARO ENRIQUE MAZA CUADRO",   "2506", 72000, 1800000),
    @("CC", "45523042", "NOHORA MILENA BERRIO CRECIAN", "2507", 80000, 2000000),
    @("CC", "73209840", "ALVARO ENRIQUE MAZA CUADRO",   "2507", 72000, 1800000),
    @("CC", "45523042", "NOHORA MILENA BERRIO CRECIAN", "2508", 80000, 2000000),
    @("CC", "73209840", "ALVARO ENRIQUE MAZA CUADRO",   "2508", 72000, 1800000)
)

$row = 16
foreach ($r in $data) {
    $ws.Cells.Item($row, 2).Value = $r[0]
    $ws.Cells.Item($row, 3).Value = $r[1]
    $ws.Cells.Item($row, 4).Value = $r[2]
    $ws.Cells.Item($row, 5).Value = $r[3]
    $ws.Cells.Item($row, 6).Value = $r[4]
    $ws.Cells.Item($row, 7).Value = $r[5]
    $row = $row + 1
}

# Update the "VALOR MORA" total to match the new period totals.
$ws.Range("E11").Value = 456000
